$wb = $excel.ActiveWorkbook

# --- pedido_cliente sheet: rename "tipo_produto" header to "produto" and
# --- insert a new "tipo" row right after "qtd_produto" (pushing
# --- valor_unitario / valor_total down by one row). ---
$ws = $wb.Worksheets.Item("pedido_cliente")

$ws.Range("A4").Value = "produto"
$ws.Rows("6").Insert()
$ws.Range("A6").Value = "tipo"

# This sheet becomes the active sheet/tab, with A9 selected.
[void]$ws.Activate()
[void]$ws.Range("A9").Select()
